$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.317.62"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.931.83"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'251.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "'0.7131"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.3264"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "'27.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").Value = "'0.07201"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").Value = "'0.08092"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "1.929.60"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "'5.429"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "'94.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'14.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "30.320.69"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'252.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.10%  "
$ws.Range("D19").Value = "'0.000008120"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").Value = "'5.791"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "2.181.61"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'9.706"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'164.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("D27").Value = "'19.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").Value = "'2.319"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("D30").Value = "'1.360"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'4.430"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "'4.206"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'0.05206"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("E35").Value = "  +4.79%  "
$ws.Range("D36").Value = "'0.7490"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "'2.764"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "'2.801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'78.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "'6.424"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "'0.4526"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").Value = "'2.028"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "'0.8411"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'101.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "'9.812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "'7.428"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "'36.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'0.06088"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("E51").Value = "  +1.10%  "
